$wb = $excel.ActiveWorkbook

# --- "PF results (NR)" sheet: add column D (angle [rad]) + formatted column E ---
$wsPF = $wb.Worksheets.Item("PF results (NR)")

# Header for new column D
$wsPF.Range("D1").Value = "angle [rad]"

# D2 is entered as its own (non-shared) formula, D3:D15 as a second fill
# producing a shared-formula group (matches how the author actually typed it).
$wsPF.Range("D2").Formula = "=C2*PI()/180"
$wsPF.Range("D3:D15").Formula = "=C3*PI()/180"

# Column E gets number formatting only (no values) - mostly format "0",
# except E3 which ended up with the workbook's existing "0.00" format.
$wsPF.Range("E2").NumberFormat = "0"
$wsPF.Range("E3").NumberFormat = "0.00"
$wsPF.Range("E4:E15").NumberFormat = "0"

# Column width tweaks: column C loses its auto bestFit sizing and both C & D
# get explicit custom widths.
$wsPF.Columns.Item(3).ColumnWidth = 10.333333333333334
$wsPF.Columns.Item(4).ColumnWidth = 8.833333333333334

# --- Selection / active-sheet bookkeeping ---
# Previously "Generators" was the active tab with C3 selected; now it is
# deselected and the selection left on E14.
$wsGen = $wb.Worksheets.Item("Generators")
$wsGen.Range("E14").Select()

# "PF results (NR)" becomes the active tab, with G15 selected.
$wsPF.Activate()
$wsPF.Range("G15").Select()
